# Trade #16 closed at 2026-02-18 00:11:07 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.73
$wsSummary.Range("B4").Value = 0.83
$wsSummary.Range("B5").Value = 0.38
$wsSummary.Range("B6").Value = 44
$wsSummary.Range("B7").Value = 25
$wsSummary.Range("B9").Value = 56.82

# --- Strategy Status sheet (MarketMaking row) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C6").Value = 99.73
$wsStatus.Range("D6").Value = 15
$wsStatus.Range("E6").Value = -0.08
$wsStatus.Range("F6").Value = -0.27
$wsStatus.Range("G6").Value = 60

# --- All Trades sheet (Trade #45 at row 46) ---
$wsAll = $wb.Worksheets.Item("All Trades")
$wsAll.Range("G46").Value = 0.44
$wsAll.Range("H46").Value = "CLOSED"
$wsAll.Range("I46").Value = 7.3171
$wsAll.Range("J46").Value = 0.03
$wsAll.Range("K46").Value = 99.73
$wsAll.Range("L46").Value = "early_exit"
$wsAll.Range("M46").Value = 0.14

# --- MarketMaking sheet (Trade #45 at row 17) ---
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("G17").Value = 0.44
$wsMM.Range("H17").Value = "CLOSED"
$wsMM.Range("I17").Value = 7.3171
$wsMM.Range("J17").Value = 0.03
$wsMM.Range("K17").Value = 99.73
$wsMM.Range("P17").Value = "early_exit"
$wsMM.Range("Q17").Value = 0.14
